$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows (old rows 6 and 7) that no longer exist in the
# updated TPM output; this shifts dimension from A1:T7 to A1:T5.
$ws.Rows("6:7").Delete()

# Refresh data rows 2-5 with the newly recomputed TPM-derived values.
$data = @(
    @("ECs",   "Fgf13", "Scn8a", "FAPs",  1, 0.3333333333333333, 0.01028766666666667, 0.030863,            0.004475330795722314, 0.004475330795722313, 3, 1, 0.654434, 1.963302,            0.7645131972140643, 0.7645131972140642, 0.006732598847333333, 0.060593389626,      0.003421449455228229, 0.003421449455228228),
    @("ECs",   "Fgf13", "Scn8a", "MuSCs", 1, 0.3333333333333333, 0.01028766666666667, 0.030863,            0.004475330795722314, 0.004475330795722313, 3, 1, 0.20158,  0.6047400000000001, 0.2354868027859358, 0.2354868027859358, 0.002073787846666667, 0.01866409062,       0.001053881340494086, 0.001053881340494085),
    @("MuSCs", "Fgf13", "Scn8a", "FAPs",  3, 1,                  2.288462333333333,  6.865386999999999,  0.9955246692042776,  0.9955246692042776,  3, 1, 0.654434, 1.963302,            0.7645131972140643, 0.7645131972140642, 1.497647558652666,   13.478828027874,     0.761091747758836,   0.7610917477588359),
    @("MuSCs", "Fgf13", "Scn8a", "MuSCs", 3, 1,                  2.288462333333333,  6.865386999999999,  0.9955246692042776,  0.9955246692042776,  3, 1, 0.20158,  0.6047400000000001, 0.2354868027859358, 0.2354868027859358, 0.4613082371533333,  4.15177413438,      0.2344329214454417,  0.2344329214454416)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
